$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Round 22 (and the round-8/round-18 catch-up matches) now have final scores
# and the "Data availability" flag for those rows flips from "N" to "Y".
# ---------------------------------------------------------------------------

# Row 122 (round 21 catch-up match) - only data-availability flips to Y
$ws.Cells.Item(122,9).Value = "Y"

# Row 123 (round 8 catch-up match) - only data-availability flips to Y
$ws.Cells.Item(123,9).Value = "Y"

# Row 124 - Gold Coast United 3-2 Eastern Suburbs
$ws.Cells.Item(124,5).Value = 3
$ws.Cells.Item(124,8).Value = 2
$ws.Cells.Item(124,9).Value = "Y"

# Row 125 - Brisbane City 4-0 Brisbane Roar Youth
$ws.Cells.Item(125,5).Value = 4
$ws.Cells.Item(125,8).Value = 0
$ws.Cells.Item(125,9).Value = "Y"

# Row 126 - Gold Coast Knights 0-2 Moreton Bay United
$ws.Cells.Item(126,5).Value = 0
$ws.Cells.Item(126,8).Value = 2
$ws.Cells.Item(126,9).Value = "Y"

# Row 127 - Olympic FC 0-4 Sunshine Coast Wanderers
$ws.Cells.Item(127,5).Value = 0
$ws.Cells.Item(127,8).Value = 4
$ws.Cells.Item(127,9).Value = "Y"

# Row 128 - Lions 1-3 Logan Lightning
$ws.Cells.Item(128,5).Value = 1
$ws.Cells.Item(128,8).Value = 3
$ws.Cells.Item(128,9).Value = "Y"

# Row 129 - Capalaba FC 2-2 Peninsula Power
$ws.Cells.Item(129,5).Value = 2
$ws.Cells.Item(129,8).Value = 2
$ws.Cells.Item(129,9).Value = "Y"

# Row 130 (round 18 catch-up) - date moves a day later, result recorded
$ws.Cells.Item(130,2).Value = "8/24/2022"
$ws.Cells.Item(130,5).Value = 1
$ws.Cells.Item(130,8).Value = 3

# ---------------------------------------------------------------------------
# Append the Finals series fixtures: two Semi-finals (SF) and the Grand
# Final (GF). Copy formatting from the last existing data row (133) first so
# the new rows keep the same style (fills/fonts/date number format), then
# fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A133:I133").Copy()
$ws.Range("A134:I136").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 134 - Semi-final: Lions v Olympic FC
$ws.Cells.Item(134,1).Value = "SF"
$ws.Cells.Item(134,2).Value = "9/3/2022"
$ws.Cells.Item(134,3).Value = "Lions"
$ws.Cells.Item(134,4).Value = "LIO"
$ws.Cells.Item(134,5).Value = -1
$ws.Cells.Item(134,6).Value = "Olympic FC"
$ws.Cells.Item(134,7).Value = "BOL"
$ws.Cells.Item(134,8).Value = -1
$ws.Cells.Item(134,9).Value = "N"

# Row 135 - Semi-final: Gold Coast Knights v Peninsula Power
$ws.Cells.Item(135,1).Value = "SF"
$ws.Cells.Item(135,2).Value = "9/4/2022"
$ws.Cells.Item(135,3).Value = "Gold Coast Knights"
$ws.Cells.Item(135,4).Value = "GCK"
$ws.Cells.Item(135,5).Value = -1
$ws.Cells.Item(135,6).Value = "Peninsula Power"
$ws.Cells.Item(135,7).Value = "PEN"
$ws.Cells.Item(135,8).Value = -1
$ws.Cells.Item(135,9).Value = "N"

# Row 136 - Grand Final: Lions v Gold Coast Knights
$ws.Cells.Item(136,1).Value = "GF"
$ws.Cells.Item(136,2).Value = "9/10/2022"
$ws.Cells.Item(136,3).Value = "Lions"
$ws.Cells.Item(136,4).Value = "LIO"
$ws.Cells.Item(136,5).Value = -1
$ws.Cells.Item(136,6).Value = "Gold Coast Knights"
$ws.Cells.Item(136,7).Value = "GCK"
$ws.Cells.Item(136,8).Value = -1
$ws.Cells.Item(136,9).Value = "N"

# ---------------------------------------------------------------------------
# View state: scroll back to the top (drop the old topLeftCell) and leave the
# selection where the author last clicked while entering the new rows.
# ---------------------------------------------------------------------------
$null = $ws.Range("R142").Select()
